$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column S carrying the 2022 figures, matching formatting of column R
# by copying the whole column R (rows 3-8) to column S first, then overwrite
# the copied values where they differ from a straight column duplication.
$ws.Range("R3:R8").Copy($ws.Range("S3"))

# Row 3 - header year
$ws.Range("S3").Value = 2022

# Row 4 - a) Number of branches per 100 000 adults
# R4 used to be a formula (=R6/R8*100000); it becomes a hard-coded value, and
# S4 gets its own hard-coded value.
$ws.Range("R4").Value = 6.9132648934880807
$ws.Range("S4").Value = 6.9031689452913012

# Row 5 - b) Number of ATMs per 100 000 adults
# R5 used to be a formula (=R7/R8*100000); it becomes a hard-coded value, and
# S5 gets its own hard-coded value.
$ws.Range("R5").Value = 42.321589572314856
$ws.Range("S5").Value = 44.306188104841333

# Row 6 - Total branches of commercial banks
$ws.Range("S6").Value = 318

# Row 7 - Total ATMs
$ws.Range("S7").Value = 2041

# Row 8 - Number of adult resident population
$ws.Range("R8").Value = 4513063
$ws.Range("S8").Value = 4606580

# Update the selected cell shown when the workbook is opened
$null = $ws.Range("R13").Select()
